$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 — swap in the Astana/Nur-Sultan monument address + coords (was row 5 content)
$ws.Range("A1").Value = "Заречный, Нур-Султан, район  Есиль, Нур-Султан, 010000, Казахстан"
$ws.Range("B1").Value = "51.12827785, 71.430515"

# Row 2 — unchanged (St. Petersburg Bogatyrsky prospect address + coords)
$ws.Range("A2").Value = "4, Богатырский проспект, Комендантский аэродром, округ Комендантский аэродром, Приморский район, Санкт-Петербург, Северо-Западный федеральный округ, 190000, РФ"
$ws.Range("B2").Value = "59.9992052, 30.2891141"

# Row 3 — now holds the old row-1 TfL Office address + coords; also gets wrap-text formatting
$ws.Range("A3").Value = "TfL Office, 208-216, Baker Street, Marylebone, City of Westminster, London, Greater London, England, NW1 5RT, UK"
$ws.Range("A3").WrapText = $true
$ws.Range("B3").Value = "51.523210, -0.157847"

# Row 4 — A4 removed entirely, only B4 remains with new text (no wrap)
$ws.Range("A4").ClearContents()
$ws.Range("B4").Value = "esrtdnjm"
$ws.Range("B4").WrapText = $false

# Row 5 — new short test values (no wrap)
$ws.Range("A5").Value = " "
$ws.Range("A5").WrapText = $false
$ws.Range("B5").Value = "ывап"
$ws.Range("B5").WrapText = $false

# Row 6 — numeric coordinate value + test string
$ws.Range("A6").Value = 51.12827785
$ws.Range("B6").Value = "?"

# Row 7 — numeric id value + test string, ends up as the active selection
$ws.Range("A7").Value = 3546575
$ws.Range("B7").Value = "½"

# Column widths shrank a bit
$ws.Columns.Item(1).ColumnWidth = 117
$ws.Columns.Item(2).ColumnWidth = 21.5

$ws.Range("B7").Select()
